$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set for rows 2-16 (driver_name, asset_id, job_site,
# scheduled_start, scheduled_end, actual_start, actual_end, status, status_reason)
$data = @(
    @("Christopher Thomas", "TRK-1011", "Downtown Construction", "06:30 AM", "03:30 PM", "06:29 AM", "03:33 PM", "On Time", ""),
    @("Matthew Harris", "TRK-1013", "North Dallas Site", "06:45 AM", "03:45 PM", "06:52 AM", "03:37 PM", "On Time", ""),
    @("Mark Thompson", "TRK-1015", "Downtown Construction", "06:00 AM", "03:00 PM", "N/A", "N/A", "On Time", ""),
    @("William Brown", "TRK-1004", "Richardson Development", "06:30 AM", "03:30 PM", "06:32 AM", "02:54 PM", "Early End", "36 minutes early"),
    @("Richard Wilson", "TRK-1007", "Downtown Construction", "06:15 AM", "03:15 PM", "06:14 AM", "03:18 PM", "On Time", ""),
    @("Anthony Martin", "TRK-1014", "West Plano Project", "07:00 AM", "04:00 PM", "07:16 AM", "03:58 PM", "Late", "16 minutes late"),
    @("James Davis", "TRK-1005", "North Dallas Site", "06:00 AM", "03:00 PM", "07:03 AM", "03:12 PM", "Not On Job", "At incorrect location: North Richland Hills"),
    @("Charles Anderson", "TRK-1010", "West Plano Project", "06:00 AM", "03:00 PM", "06:05 AM", "02:48 PM", "On Time", ""),
    @("Michael Johnson", "TRK-1002", "West Plano Project", "06:00 AM", "03:00 PM", "06:18 AM", "03:07 PM", "Late", "18 minutes late"),
    @("Robert Williams", "TRK-1003", "Downtown Construction", "07:00 AM", "04:00 PM", "07:28 AM", "03:45 PM", "Late", "28 minutes late"),
    @("David Miller", "TRK-1006", "West Plano Project", "06:45 AM", "03:45 PM", "06:47 AM", "03:42 PM", "On Time", ""),
    @("Joseph Moore", "TRK-1008", "Richardson Development", "06:30 AM", "03:30 PM", "06:55 AM", "03:25 PM", "Late", "25 minutes late"),
    @("Thomas Taylor", "TRK-1009", "North Dallas Site", "07:00 AM", "04:00 PM", "06:43 AM", "03:52 PM", "On Time", ""),
    @("Daniel Jackson", "TRK-1012", "Richardson Development", "06:15 AM", "03:15 PM", "06:16 AM", "03:12 PM", "On Time", ""),
    @("John Smith", "TRK-1001", "North Dallas Site", "06:30 AM", "03:30 PM", "06:25 AM", "03:35 PM", "On Time", "")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $val = $rowValues[$col - 1]
        if ($val -ne "") {
            $ws.Cells.Item($rowNum, $col).Value = $val
        }
    }
}
